$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the extra data that was populated in row 4 (columns I:N) —
# these cells (amounts + address text) are removed in the committed revision.
$ws.Range("I4:N4").Clear()
